# Refresh the cached statistics that are sourced from the external workbook
# "Create/_Test_Suite_Statistics_for_folders.xlsx" ([1]Sheet1!...).
#
# The external workbook isn't reachable from this sandbox (no real linked
# file to pull fresh numbers from), so we push the refreshed figures
# straight onto the cells that consume them. Everything else on Sheet1
# (H1, H2, H3, H5, H6, H7, J2, ...) is a local formula (SUM/ratio) that
# depends on these cells, so it recalculates on its own once the inputs
# change - no need to touch those formulas directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values pulled from the refreshed "Create" external source
# (was 19 -> 15, 31 -> 25, 149 -> 120, 230 -> 194 for H3/H1/H6/H5 respectively).
$ws.Range("C2").Value = 15    # was: =[1]Sheet1!$H$3  (19 -> 15)
$ws.Range("D2").Value = 25    # was: =[1]Sheet1!$H$1  (31 -> 25)
$ws.Range("E2").Value = 120   # was: =[1]Sheet1!$H$6  (149 -> 120)
$ws.Range("F2").Value = 194   # was: =[1]Sheet1!$H$5  (230 -> 194)

$wb.Save()
